$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so numeric-looking
# strings (e.g. "0.629") are not silently coerced to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.201.40"
$ws.Range("E2").Value = "  -5.88%  "
$ws.Range("D3").Value = "2.222.40"
$ws.Range("E3").Value = "  -5.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "245.86"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  -5.86%  "
$ws.Range("D7").Value = "70.78"
$ws.Range("E7").Value = "  -4.72%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  -6.99%  "
$ws.Range("D10").Value = "38.75"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "58.45"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0950"
$ws.Range("E12").Value = "  -7.06%  "
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("D14").Value = "6.76"
$ws.Range("E14").Value = "  -7.87%  "
$ws.Range("D15").Value = "2.552.60"
$ws.Range("E15").Value = "  -5.46%  "
$ws.Range("E16").Value = "  -9.35%  "
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -8.54%  "
$ws.Range("D18").Value = "2.221.63"
$ws.Range("E18").Value = "  -5.33%  "
$ws.Range("D19").Value = "41.279.85"
$ws.Range("E19").Value = "  -5.44%  "
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  -7.77%  "
$ws.Range("D21").Value = "72.57"
$ws.Range("E21").Value = "  -5.76%  "
$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  -7.47%  "
$ws.Range("D23").Value = "232.32"
$ws.Range("E23").Value = "  -8.22%  "
$ws.Range("E24").Value = "  +11.74%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "3.69"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "2.43"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  -7.38%  "
$ws.Range("D30").Value = "171.33"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "20.56"
$ws.Range("E31").Value = "  -7.76%  "
$ws.Range("E32").Value = "  -7.75%  "
$ws.Range("E33").Value = "  -6.96%  "
$ws.Range("D34").Value = "0.0715"
$ws.Range("E34").Value = "  -5.22%  "
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("D36").Value = "4.62"
$ws.Range("E36").Value = "  -10.08%  "
$ws.Range("D37").Value = "3.91"
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("D38").Value = "24.14"
$ws.Range("E38").Value = "  +14.53%  "
$ws.Range("D39").Value = "0.0277"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("E41").Value = "  -11.12%  "
$ws.Range("D42").Value = "65.57"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("E43").Value = "  -8.72%  "
$ws.Range("D44").Value = "0.203"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "8.79"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "10.95"
$ws.Range("E46").Value = "  +11.92%  "
$ws.Range("E47").Value = "  -6.79%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "4.53"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("E50").Value = "  -5.80%  "
$ws.Range("E51").Value = "  -4.97%  "

# Restore default (unstyled) formatting now that values are committed as text.
$ws.Range("D2:E51").Style = "Normal"
